$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 19
$ws.Cells.Item($row, 1).Value2 = 42601.898912037039
$ws.Cells.Item($row, 2).Value = "Named"
$ws.Cells.Item($row, 3).Value2 = 12991
$ws.Cells.Item($row, 4).Value2 = 7780
$ws.Cells.Item($row, 5).Value2 = 483
$ws.Cells.Item($row, 6).Value2 = 65
$ws.Cells.Item($row, 7).Value2 = 38
$ws.Cells.Item($row, 8).Value2 = 63
$ws.Cells.Item($row, 9).Value2 = 36
$ws.Cells.Item($row, 10).Value2 = 0
$ws.Cells.Item($row, 11).Value2 = 0
$ws.Cells.Item($row, 12).Value2 = 0
$ws.Cells.Item($row, 13).Value2 = 0
